$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.944.65'
$ws.Range("E2").Value = '  -3.92%  '
$ws.Range("D3").Value = '1.636.98'
$ws.Range("E3").Value = '  -6.24%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4754'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2577'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06080'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07008'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("D11").Value = '1.641.64'
$ws.Range("E11").Value = '  -5.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.66'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5993'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.343'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '73.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9990'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9979'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '24.943.64'
$ws.Range("E18").Value = '  -3.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006568'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.95%  '
$ws.Range("D21").Value = '1.850.31'
$ws.Range("E21").Value = '  -6.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.356'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.546'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.238'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '133.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.383'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '103.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.634'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.888'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07703'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.543'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9983'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04308'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.582'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9229'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5807'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.527'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01528'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9978'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8152'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.761'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3686'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.699'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1088'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05185'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.046'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9986'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9955'
$ws.Range("D51").Style = "Normal"
